$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new value, derived from the authoritative diff.
$changes = @{
    "B7"   = 5.293599999999998
    "A8"   = -22.36890000000002
    "A10"  = -21.64289999999999
    "A12"  = -21.6016
    "B14"  = 5.864600000000003
    "B15"  = 4.748699999999992
    "A18"  = -21.83370000000002
    "B18"  = 6.248099999999998
    "B20"  = 8.642700000000003
    "A25"  = -21.57249999999999
    "B29"  = 5.042400000000004
    "B30"  = 5.086000000000001
    "B31"  = 5.217400000000005
    "B35"  = 8.460300000000005
    "A37"  = -18.9653
    "B40"  = 9.142399999999995
    "B44"  = 5.012800000000004
    "B50"  = 4.620899999999998
    "B54"  = 4.506099999999997
    "A55"  = -22.0843
    "A68"  = -21.47430000000002
    "B68"  = 4.452499999999999
    "B76"  = 5.952199999999997
    "A77"  = -20.11049999999999
    "A78"  = -20.07179999999998
    "A79"  = -19.95489999999999
    "A80"  = -19.67249999999999
    "A81"  = -21.70590000000001
    "A82"  = -21.8491
    "A84"  = -22.0673
    "B87"  = 4.738899999999994
    "B88"  = 4.512299999999997
    "B92"  = 4.622799999999997
    "B96"  = 5.002300000000008
    "B98"  = 6.478499999999999
    "A101" = -21.66699999999999
    "B101" = 5.550900000000002
    "A102" = -21.72919999999997
    "B102" = 5.434600000000005
}

foreach ($ref in $changes.Keys) {
    $ws.Range($ref).Value = $changes[$ref]
}
